# Generate Report for Handback
#
# Replaces the two source-file GUID identifiers (and their derived hash
# filenames / handoff-handback timestamps) across the Overview, zh-cn and
# de-de sheets, updating both the cell values and the visible text of the
# corresponding hyperlinks (their Address/r:id stay untouched).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New values
# ---------------------------------------------------------------------
$newId1 = "0e71f548-af60-45c4-b839-55682e57a859"
$newId2 = "ffff346c3780-2602-4cf0-9d3f-bc4d5f73ec6c"

$newMd1 = "$newId1.md"
$newMd2 = "$newId2.md"

$newXlfZhCn = "$newId1.0aa05150851e84d6a81556a899c19be5e8830348.zh-cn.xlf"
$newXlfDeDe = "$newId1.0aa05150851e84d6a81556a899c19be5e8830348.de-de.xlf"

$zhHandoffTime = "2016-03-11 09:46:44"
$zhHandbackTime = "2016-03-11 09:47:01"
$deHandoffTime = "2016-03-11 09:46:47"
$deHandbackTime = "2016-03-11 09:47:07"

function Get-AbsAddress([string]$CellRef) {
    if ($CellRef -match '^([A-Za-z]+)(\d+)$') {
        return '$' + $matches[1] + '$' + $matches[2]
    }
    return $null
}

function Set-CellAndLink {
    param(
        $Sheet,
        [string]$CellRef,
        [string]$NewValue
    )
    $Sheet.Range($CellRef).Value = $NewValue
    $target = Get-AbsAddress $CellRef
    foreach ($hl in $Sheet.Hyperlinks) {
        if ($hl.Range.Address() -eq $target) {
            $hl.TextToDisplay = $NewValue
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
Set-CellAndLink $wsOverview "A2" $newMd1
Set-CellAndLink $wsOverview "A3" $newMd2

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
Set-CellAndLink $wsZh "A2" $newMd1
Set-CellAndLink $wsZh "F2" $newMd1
Set-CellAndLink $wsZh "D2" $newXlfZhCn
Set-CellAndLink $wsZh "G2" $newXlfZhCn
$wsZh.Range("E2").Value = $zhHandoffTime
$wsZh.Range("H2").Value = $zhHandbackTime

Set-CellAndLink $wsZh "A3" $newMd2
Set-CellAndLink $wsZh "F3" $newMd2
Set-CellAndLink $wsZh "D3" $newXlfZhCn
Set-CellAndLink $wsZh "G3" $newXlfZhCn
$wsZh.Range("E3").Value = $zhHandoffTime
$wsZh.Range("H3").Value = $zhHandbackTime

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
Set-CellAndLink $wsDe "A2" $newMd1
Set-CellAndLink $wsDe "F2" $newMd1
Set-CellAndLink $wsDe "D2" $newXlfDeDe
Set-CellAndLink $wsDe "G2" $newXlfDeDe
$wsDe.Range("E2").Value = $deHandoffTime
$wsDe.Range("H2").Value = $deHandbackTime

Set-CellAndLink $wsDe "A3" $newMd2
Set-CellAndLink $wsDe "F3" $newMd2
Set-CellAndLink $wsDe "D3" $newXlfDeDe
Set-CellAndLink $wsDe "G3" $newXlfDeDe
$wsDe.Range("E3").Value = $deHandoffTime
$wsDe.Range("H3").Value = $deHandbackTime
